$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47, shifting existing rows 47:80 down to 48:81.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new data record.
$ws.Range("A47").Value = 3
$ws.Range("B47").Value = "Femacal de La Calera"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 45216
$ws.Range("E47").Value = 5
$ws.Range("F47").Value = 300000000
$ws.Range("G47").Value = "Espárragos"
$ws.Range("H47").Value = "Verde"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 1600
$ws.Range("K47").Value = 1800
$ws.Range("L47").Value = 1800
$ws.Range("M47").Value = 1800
$ws.Range("N47").Value = "$/kilo"
$ws.Range("O47").Value = "Provincia de Linares"
$ws.Range("P47").Value = 1800
$ws.Range("Q47").Value = 1
$ws.Range("R47").Value = "Hortaliza"

# Match the date style/number format used by the other rows in column D.
$ws.Range("D47").NumberFormat = $ws.Range("D48").NumberFormat
